$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.512729
$ws.Range("H2").Value = 7.538187
$ws.Range("I2").Value = 0.02190726325199687
$ws.Range("J2").Value = 0.02190726325199687
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 59.39035423516566
$ws.Range("R2").Value = 534.5131881164909
$ws.Range("S2").Value = 0.001495793888247391
$ws.Range("T2").Value = 0.001495793888247391

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.512729
$ws.Range("H3").Value = 7.538187
$ws.Range("I3").Value = 0.02190726325199687
$ws.Range("J3").Value = 0.02190726325199687
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 455.5285982666153
$ws.Range("R3").Value = 4099.757384399538
$ws.Range("S3").Value = 0.01147285450615571
$ws.Range("T3").Value = 0.01147285450615571

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.512729
$ws.Range("H4").Value = 7.538187
$ws.Range("I4").Value = 0.02190726325199687
$ws.Range("J4").Value = 0.02190726325199687
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 279.1655646876103
$ws.Range("R4").Value = 2512.490082188493
$ws.Range("S4").Value = 0.007031009510659039
$ws.Range("T4").Value = 0.00703100951065904

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.512729
$ws.Range("H5").Value = 7.538187
$ws.Range("I5").Value = 0.02190726325199687
$ws.Range("J5").Value = 0.02190726325199687
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 75.74128908100732
$ws.Range("R5").Value = 681.6716017290659
$ws.Range("S5").Value = 0.001907605346934728
$ws.Range("T5").Value = 0.001907605346934728

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 92.89399466666667
$ws.Range("H6").Value = 278.681984
$ws.Range("I6").Value = 0.8098976036382196
$ws.Range("J6").Value = 0.8098976036382197
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 2195.623662389746
$ws.Range("R6").Value = 19760.61296150771
$ws.Range("S6").Value = 0.05529854969528578
$ws.Range("T6").Value = 0.05529854969528578

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 92.89399466666667
$ws.Range("H7").Value = 278.681984
$ws.Range("I7").Value = 0.8098976036382196
$ws.Range("J7").Value = 0.8098976036382197
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 16840.60285764725
$ws.Range("R7").Value = 151565.4257188252
$ws.Range("S7").Value = 0.4241441418100684
$ws.Range("T7").Value = 0.4241441418100684

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 92.89399466666667
$ws.Range("H8").Value = 278.681984
$ws.Range("I8").Value = 0.8098976036382196
$ws.Range("J8").Value = 0.8098976036382197
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 10320.57355855242
$ws.Range("R8").Value = 92885.16202697177
$ws.Range("S8").Value = 0.259931954454477
$ws.Range("T8").Value = 0.259931954454477

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 92.89399466666667
$ws.Range("H9").Value = 278.681984
$ws.Range("I9").Value = 0.8098976036382196
$ws.Range("J9").Value = 0.8098976036382197
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 2800.10733506779
$ws.Range("R9").Value = 25200.96601561011
$ws.Range("S9").Value = 0.0705229576783885
$ws.Range("T9").Value = 0.07052295767838851

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.610639333333333
$ws.Range("H10").Value = 4.831918
$ws.Range("I10").Value = 0.0140423817607685
$ws.Range("J10").Value = 0.0140423817607685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 38.06874539664155
$ws.Range("R10").Value = 342.618708569774
$ws.Range("S10").Value = 0.000958792003025735
$ws.Range("T10").Value = 0.0009587920030257351

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.610639333333333
$ws.Range("H11").Value = 4.831918
$ws.Range("I11").Value = 0.0140423817607685
$ws.Range("J11").Value = 0.0140423817607685
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 291.9902137581925
$ws.Range("R11").Value = 2627.911923823732
$ws.Range("S11").Value = 0.007354008622985195
$ws.Range("T11").Value = 0.007354008622985194

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.610639333333333
$ws.Range("H12").Value = 4.831918
$ws.Range("I12").Value = 0.0140423817607685
$ws.Range("J12").Value = 0.0140423817607685
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 178.9429098792891
$ws.Range("R12").Value = 1610.486188913602
$ws.Range("S12").Value = 0.004506821257249868
$ws.Range("T12").Value = 0.004506821257249868

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.610639333333333
$ws.Range("H13").Value = 4.831918
$ws.Range("I13").Value = 0.0140423817607685
$ws.Range("J13").Value = 0.0140423817607685
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 48.54956477648044
$ws.Range("R13").Value = 436.946082988324
$ws.Range("S13").Value = 0.001222759877507703
$ws.Range("T13").Value = 0.001222759877507703

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.68108066666667
$ws.Range("H14").Value = 53.04324200000001
$ws.Range("I14").Value = 0.154152751349015
$ws.Range("J14").Value = 0.154152751349015
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 417.9064451653452
$ws.Range("R14").Value = 3761.158006488106
$ws.Range("S14").Value = 0.01052531028965285
$ws.Range("T14").Value = 0.01052531028965285

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.68108066666667
$ws.Range("H15").Value = 53.04324200000001
$ws.Range("I15").Value = 0.154152751349015
$ws.Range("J15").Value = 0.154152751349015
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 3205.374671094902
$ws.Range("R15").Value = 28848.37203985411
$ws.Range("S15").Value = 0.08072994182829479
$ws.Range("T15").Value = 0.08072994182829478

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.68108066666667
$ws.Range("H16").Value = 53.04324200000001
$ws.Range("I16").Value = 0.154152751349015
$ws.Range("J16").Value = 0.154152751349015
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 1964.377721830404
$ws.Range("R16").Value = 17679.39949647364
$ws.Range("S16").Value = 0.04947443449972641
$ws.Range("T16").Value = 0.04947443449972641

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.68108066666667
$ws.Range("H17").Value = 53.04324200000001
$ws.Range("I17").Value = 0.154152751349015
$ws.Range("J17").Value = 0.154152751349015
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 532.961509991173
$ws.Range("R17").Value = 4796.653589920556
$ws.Range("S17").Value = 0.01342306473134094
$ws.Range("T17").Value = 0.01342306473134094
